$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Restore header row (this reverts the previous edit that replaced the
# "Username"/"Password" header with an extra "Admin"/"Paswword123" row).
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Keep the existing credential rows intact (unchanged by this revert).
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"

$ws.Range("A3").Value = "Shubham"
$ws.Range("B3").Value = "shubham123"

$ws.Range("A4").Value = "Abhijit"
$ws.Range("B4").Value = "Abhijit123"

$ws.Range("A5").Value = "Admin"
$ws.Range("B5").Value = "Admin1234"

# Restore the originally selected cell.
$ws.Range("B5").Select()
